$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1635
$ws.Range("I33").Value = 1607
$ws.Range("J33").Value = 1700.3334
$ws.Range("K33").Value = 1607
$ws.Range("L33").Value = 1700.3334
$ws.Range("M33").Value = -1378
$ws.Range("N33").Value = -2158.3334

$ws.Range("H48").Value = 1669.5
$ws.Range("I48").Value = 1008.5
$ws.Range("K48").Value = 3025.5
$ws.Range("M48").Value = -2733.5

$ws.Range("H56").Value = 1669.5
$ws.Range("I56").Value = 1008.5
$ws.Range("K56").Value = 3025.5
$ws.Range("M56").Value = -2491.5

$ws.Range("H69").Value = 44682.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 44682.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 134047.5
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -135795.5

$ws.Range("H72").Value = 44682.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 44682.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 402142.5
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -410878.5

$ws.Range("H92").Value = 854.625
$ws.Range("I92").Value = 560.41174
$ws.Range("K92").Value = 560.41174
$ws.Range("M92").Value = 687.58826

$ws.Range("H101").Value = 2849.5
$ws.Range("I101").Value = 3019.4
$ws.Range("J101").Value = 2000
$ws.Range("K101").Value = 9058.200000000001
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = -7436.200000000001
$ws.Range("N101").Value = -9244

$ws.Range("H137").Value = 3538
$ws.Range("I137").Value = 901.8125
$ws.Range("J137").Value = 6174.1875
$ws.Range("K137").Value = 2705.4375
$ws.Range("L137").Value = 18522.5625
$ws.Range("M137").Value = -155.4375
$ws.Range("N137").Value = -23622.5625

$ws.Range("H138").Value = 3265.6086
$ws.Range("J138").Value = 3261.2632
$ws.Range("L138").Value = 9783.7896
$ws.Range("N138").Value = -20063.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1831.3214
$ws.Range("I45").Value = 962.9091
$ws.Range("K45").Value = 962.9091
$ws.Range("M45").Value = -585.9091

$ws.Range("H52").Value = 91496.336
$ws.Range("J52").Value = 91496.336
$ws.Range("L52").Value = 91496.336
$ws.Range("N52").Value = -92132.336

$ws.Range("H61").Value = 16709946
$ws.Range("I61").Value = 38467184
$ws.Range("K61").Value = 38467184
$ws.Range("M61").Value = -38466972

$ws.Range("H63").Value = 5997.35
$ws.Range("I63").Value = 2867.9092
$ws.Range("J63").Value = 9822.223
$ws.Range("K63").Value = 2867.9092
$ws.Range("L63").Value = 9822.223
$ws.Range("M63").Value = -2181.9092
$ws.Range("N63").Value = -11194.223

$ws.Range("H66").Value = 5997.35
$ws.Range("I66").Value = 2867.9092
$ws.Range("J66").Value = 9822.223
$ws.Range("K66").Value = 14339.546
$ws.Range("L66").Value = 49111.115
$ws.Range("M66").Value = -10907.546
$ws.Range("N66").Value = -55975.115

$ws.Range("H74").Value = 7359167
$ws.Range("I74").Value = 9260426
$ws.Range("K74").Value = 9260426
$ws.Range("M74").Value = -9259552

$ws.Range("H77").Value = 7359167
$ws.Range("I77").Value = 9260426
$ws.Range("K77").Value = 46302130
$ws.Range("M77").Value = -46297762

$ws.Range("H136").Value = 16709946
$ws.Range("I136").Value = 38467184
$ws.Range("K136").Value = 115401552
$ws.Range("M136").Value = -115399002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 11515.667
$ws.Range("I99").Value = 18682.334
$ws.Range("J99").Value = 4349
$ws.Range("K99").Value = 18682.334
$ws.Range("L99").Value = 4349
$ws.Range("M99").Value = -17184.334
$ws.Range("N99").Value = -7345

$ws.Range("H123").Value = 135000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 135000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 135000
$ws.Range("M123").Value = ""
$ws.Range("N123").Value = -144800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""

$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -797
$ws.Range("N58").Value = ""

$ws.Range("H107").Value = 1120.1818
$ws.Range("I107").Value = 799.75
$ws.Range("K107").Value = 799.75
$ws.Range("M107").Value = 1120.25

$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1600.7142
$ws.Range("I7").Value = 1600.7142
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4802.142599999999
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4690.142599999999
$ws.Range("N7").Value = ""

$ws.Range("H120").Value = 12197.25
$ws.Range("I120").Value = 894.5
$ws.Range("K120").Value = 2683.5
$ws.Range("M120").Value = 2154.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 66669120
$ws.Range("I132").Value = 80002410
$ws.Range("K132").Value = 240007230
$ws.Range("M132").Value = -240004700

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 115325
$ws.Range("I7").Value = 3981.6
$ws.Range("K7").Value = 3981.6
$ws.Range("M7").Value = -3869.6

$ws.Range("H16").Value = 1376
$ws.Range("I16").Value = 1362.6666
$ws.Range("K16").Value = 1362.6666
$ws.Range("M16").Value = -1192.6666

$ws.Range("H40").Value = 3984.0527
$ws.Range("I40").Value = 3069.5
$ws.Range("K40").Value = 3069.5
$ws.Range("M40").Value = -2933.5

$ws.Range("H46").Value = 2318.1
$ws.Range("I46").Value = 2312.4285
$ws.Range("J46").Value = 2331.3333
$ws.Range("K46").Value = 2312.4285
$ws.Range("L46").Value = 2331.3333
$ws.Range("M46").Value = -2124.4285
$ws.Range("N46").Value = -2707.3333

$ws.Range("H55").Value = 71429460
$ws.Range("I55").Value = 100000920
$ws.Range("J55").Value = 787
$ws.Range("K55").Value = 100000920
$ws.Range("L55").Value = 787
$ws.Range("M55").Value = -100000747
$ws.Range("N55").Value = -1133

$ws.Range("H82").Value = 2279.3333
$ws.Range("J82").Value = 2980
$ws.Range("L82").Value = 2980
$ws.Range("N82").Value = -3702

$ws.Range("H85").Value = 2279.3333
$ws.Range("J85").Value = 2980
$ws.Range("L85").Value = 2980
$ws.Range("N85").Value = -5476

$ws.Range("H126").Value = 115325
$ws.Range("I126").Value = 3981.6
$ws.Range("K126").Value = 11944.8
$ws.Range("M126").Value = -9474.799999999999

$ws.Range("H132").Value = 6540864
$ws.Range("I132").Value = 792908.9
$ws.Range("K132").Value = 2378726.7
$ws.Range("M132").Value = -2376196.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 35644.75
$ws.Range("I81").Value = 24637.445
$ws.Range("J81").Value = 68666.664
$ws.Range("K81").Value = 49274.89
$ws.Range("L81").Value = 137333.328
$ws.Range("M81").Value = -48213.89
$ws.Range("N81").Value = -139455.328

$ws.Range("H84").Value = 35644.75
$ws.Range("I84").Value = 24637.445
$ws.Range("J84").Value = 68666.664
$ws.Range("K84").Value = 246374.45
$ws.Range("L84").Value = 686666.64
$ws.Range("M84").Value = -241070.45
$ws.Range("N84").Value = -697274.64

$ws.Range("H100").Value = 1637.3704
$ws.Range("I100").Value = 1656.26
$ws.Range("K100").Value = 3312.52
$ws.Range("M100").Value = -2771.52

$ws.Range("H113").Value = 606.025
$ws.Range("I113").Value = 574.3570999999999
$ws.Range("J113").Value = 679.9167
$ws.Range("K113").Value = 1723.0713
$ws.Range("L113").Value = 2039.7501
$ws.Range("M113").Value = 446.9287000000002
$ws.Range("N113").Value = -6379.7501

$ws.Range("H132").Value = 505955
$ws.Range("I132").Value = 4806.4
$ws.Range("J132").Value = 2009400.8
$ws.Range("K132").Value = 14419.2
$ws.Range("L132").Value = 6028202.4
$ws.Range("M132").Value = -11889.2
$ws.Range("N132").Value = -6033262.4

$ws.Range("H136").Value = 3472.3333
$ws.Range("I136").Value = 1458.5
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 4375.5
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -1825.5
